$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "['анализ', 'microsoft', 'управление', 'данные', 'финансовый', 'sql', 'грамотность', 'проведение', 'bi', 'мышление', 'знание', 'ведение', 'разработка', 'работа', 'аналитический', 'маркетинговый', 'python', 'исследование', 'отчетность', 'инвестиционный']"
$ws.Range("D3").Value = "['управление', 'анализ', 'microsoft', 'данные', 'грамотность', 'sql', 'работа', 'разработка', 'мышление', 'аналитический', 'bi', 'информационный', 'технический', 'бизнес', 'процесс', 'формирование', 'задание', 'excel', 'проект', 'bpmn']"
$ws.Range("D4").Value = "['управление', 'бизнес', 'разработка', 'процесс', 'bi', 'организация', 'microsoft', 'bpmn', 'анализ', 'работа', 'грамотность', 'проект', 'мышление', 'аналитический', 'система', 'технический', 'оптимизация', 'sql', 'формирование', 'задание']"
$ws.Range("D5").Value = "['проведение', 'тренинг', 'коучинг', 'продажа', 'обучаемость', 'адаптивность', 'грамотность', 'персонал', 'управление', 'коммуникативный', 'публичный', 'презентация', 'выступление', 'учебный', 'лидерство', 'организаторский', 'навык', 'потребитель', 'процесс', 'организация']"
$ws.Range("D6").Value = "['управление', 'продукт', 'грамотность', 'анализ', 'проект', 'коммуникативный', 'работа', 'проведение', 'продажа', 'мышление', 'знание', 'лидерство', 'персонал', 'разработка', 'аналитический', 'рынок', 'исследование', 'маркетинговый', 'язык', 'английский']"
$ws.Range("D7").Value = "['ведение', 'грамотность', 'документооборот', 'базовый', 'знание', 'работа', 'пк', 'коммуникативный', 'информационный', 'microsoft', 'учет', 'команда', 'управление', 'данные', 'база', 'инструментальный', 'складской', 'оргтехника', 'оформление', 'прием']"
$ws.Range("D8").Value = "['javascript', 'разработка', 'scala', 'kotlin', 'java', 'git', 'php', 'sql', 'html', 'css', 'api', 'postgresql', 'python', 'управление', 'мобильный', 'программирование', 'объектно', 'docker', 'знание', 'ориентированный']"
$ws.Range("D9").Value = "['управление', 'api', 'sql', 'bpmn', 'анализ', 'uml', 'разработка', 'системный', 'проект', 'работа', 'технический', 'формирование', 'задание', 'система', 'javascript', 'принцип', 'знание', 'данные', 'sap', 'базовый']"
$ws.Range("D10").Value = "['ведение', 'финансовый', 'отчетность', 'анализ', 'управление', 'бюджетирование', 'учет', 'грамотность', 'управленческий', 'финансы', 'microsoft', 'бухгалтерский', 'коммуникативный', 'знание', 'продажа', 'работа', 'excel', 'информационный', 'мышление', 'аналитический']"
